# Convert the field { m:'doc.html'.fromHTMLURI() } from a real Word field
# (fldChar begin/end + instrText runs) into plain literal text runs using
# "{" and "}" delimiters instead, as done by TokenIteratorFieldRewriterSplit.
#
# Strategy: locate the paragraph that holds the field, then overwrite its
# Range with a hand-built WordOpenXML fragment containing the equivalent
# plain-text runs (preserving the bookmark and run rsid attributes).

$d = $word.ActiveDocument

# Find the paragraph containing the field (begin/end fldChar pair).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range

    $newParagraphXml = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
        '<w:r><w:t>{</w:t></w:r>' +
        '<w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r>' +
        '<w:r w:rsidR="002033E1"><w:t>:</w:t></w:r>' +
        '<w:r w:rsidR="008B76C9"><w:t>''</w:t></w:r>' +
        '<w:r w:rsidR="00E806A4"><w:t>doc.html</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r w:rsidR="008B76C9"><w:t>''.fromHTMLURI()</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
        '</w:p>'

    $package = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $newParagraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    $r.InsertXML($package)
}
